$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 16: continuation of the GPTQ work session on 2025-12-08 ---
$ws.Cells.Item(16, 2).Value = 45999
$ws.Cells.Item(16, 3).Value = 0.51041666666666663
$ws.Cells.Item(16, 4).Value = 0.63541666666666663
$ws.Cells.Item(16, 6).Value = "GPTQ"

# Extend the shared "duration" formula down into the new row, keeping the
# original time number format intact.
$ws.Range("E16").NumberFormat = "[$-F400]h:mm:ss\ AM/PM"
$ws.Range("E15:E16").Formula = "=D15-C15"

# --- Row 17: new activity "Evaluate Models" ---
$ws.Cells.Item(17, 2).Value = 45999
$ws.Cells.Item(17, 6).Value = "Evaluate Models"

$ws.Range("H12").Select()
